$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the report. It goes in right above
# the current row 176 (Kiwi / Especial, Región de O'Higgins), so shift the
# existing data (rows 176-229) down by one row first.
$ws.Rows.Item(176).Insert()

# Fill in the new row 176 with the new week's data. All the descriptive
# columns (market, region, product hierarchy, unit, origin...) are the same
# as the row that used to be here (now row 177); only the date and the
# volume/price/unit-price figures differ.
$ws.Range("A176").Value = 8
$ws.Range("B176").Value = "Terminal La Palmera de La Serena"
$ws.Range("C176").Value = "Coquimbo"
$ws.Range("D176").Value = 44468
$ws.Range("E176").Value = 4
$ws.Range("F176").Value = "Fruta"
$ws.Range("G176").Value = 100101
$ws.Range("H176").Value = "Berries"
$ws.Range("I176").Value = 100101007
$ws.Range("J176").Value = "Kiwi"
$ws.Range("K176").Value = "Hayward"
$ws.Range("L176").Value = "Especial"
$ws.Range("M176").Value = 26
$ws.Range("N176").Value = 410000
$ws.Range("O176").Value = 420000
$ws.Range("P176").Value = 415000
$ws.Range("Q176").Value = "`$/bins (450 kilos)"
$ws.Range("R176").Value = "Región de O'Higgins"
$ws.Range("S176").Value = 922
$ws.Range("T176").Value = 450
